$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Delivery Date value in D2 (date serial 46481 -> 45753)
$ws.Range("D2").Value = 45753

# Update the active cell selection shown in the sheet view (D3 -> D5)
$ws.Range("D5").Select()
